$wb = $excel.ActiveWorkbook

# Rename the "fact" sheet to "facts"
$ws = $wb.Worksheets.Item("fact")
$ws.Name = "facts"

# Add a new (empty) cell at B12 with a right-aligned style, extending the
# sheet's used range from A1:D11 to A1:D12
$cell = $ws.Range("B12")
$cell.Value = $null
$cell.HorizontalAlignment = -4152   # xlRight

# Make B12 the active selection on the "facts" sheet
$cell.Select() | Out-Null
